$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1424918639367689
$ws.Range("D2").Value = 0.01666762993355775
$ws.Range("E2").Value = 0.4245806148541504
$ws.Range("F2").Value = 0.4134880438093163
$ws.Range("G2").Value = 0.2638777770573171
$ws.Range("H2").Value = 0.4207054203120748
$ws.Range("K2").Value = 0.6674196183508059
$ws.Range("N2").Value = 1.114047915013927
$ws.Range("O2").Value = 1.289226636290863
$ws.Range("B3").Value = 0.1329167770221034
$ws.Range("D3").Value = 0.01465201566882968
$ws.Range("E3").Value = 0.3704154376978721
$ws.Range("F3").Value = 0.4077945561905594
$ws.Range("G3").Value = 0.2594061449338696
$ws.Range("H3").Value = 0.4222210729690019
$ws.Range("K3").Value = 0.5860841411203523
$ws.Range("N3").Value = 1.122577175422137
$ws.Range("O3").Value = 1.282656538331864
$ws.Range("B4").Value = 0.1271055766230376
$ws.Range("D4").Value = 0.01340807465636118
$ws.Range("E4").Value = 0.3372382214265031
$ws.Range("F4").Value = 0.4046176420065564
$ws.Range("G4").Value = 0.2569064568466288
$ws.Range("H4").Value = 0.4233858122323966
$ws.Range("K4").Value = 0.5359062117629207
$ws.Range("N4").Value = 1.1282589357684
$ws.Range("O4").Value = 1.279624698882856
$ws.Range("B5").Value = 0.1247547377220002
$ws.Range("D5").Value = 0.01289959950418762
$ws.Range("E5").Value = 0.3237369150036642
$ws.Range("F5").Value = 0.4034030728197351
$ws.Range("G5").Value = 0.2559494082885081
$ws.Range("H5").Value = 0.4239192941078826
$ws.Range("K5").Value = 0.5153996128085794
$ws.Range("N5").Value = 1.130686268287917
$ws.Range("O5").Value = 1.278640633811051
$ws.Range("B6").Value = 0.12436543145337
$ws.Range("D6").Value = 0.0128150745055251
$ws.Range("E6").Value = 0.3214961165082713
$ws.Range("F6").Value = 0.4032062251180832
$ws.Range("G6").Value = 0.2557942041936414
$ws.Range("H6").Value = 0.4240114313411283
$ws.Range("K6").Value = 0.5119909946703842
$ws.Range("N6").Value = 1.131096091560934
$ws.Range("O6").Value = 1.278492399506376
$ws.Range("B7").Value = 0.1270738022024887
$ws.Range("D7").Value = 0.0134012234468841
$ws.Range("E7").Value = 0.3370560644852389
$ws.Range("F7").Value = 0.4046009379772357
$ws.Range("G7").Value = 0.2568933006824636
$ws.Range("H7").Value = 0.4233927687521799
$ws.Range("K7").Value = 0.5356298887964499
$ws.Range("N7").Value = 1.128291218091356
$ws.Range("O7").Value = 1.279610410169283
$ws.Range("B8").Value = 0.1391763822185084
$ws.Range("D8").Value = 0.01597398452878451
$ws.Range("E8").Value = 0.405886824367073
$ws.Range("F8").Value = 0.4114586500460931
$ws.Range("G8").Value = 0.2622847714219105
$ws.Range("H8").Value = 0.4211794058253417
$ws.Range("K8").Value = 0.6394250778587889
$ws.Range("N8").Value = 1.116896571435248
$ws.Range("O8").Value = 1.286752896241723
$ws.Range("B9").Value = 0.1634416258150964
$ws.Range("D9").Value = 0.02096743215687269
$ws.Range("E9").Value = 0.5415777296337865
$ws.Range("F9").Value = 0.427444864918499
$ws.Range("G9").Value = 0.2748205230819138
$ws.Range("H9").Value = 0.4186984486297121
$ws.Range("K9").Value = 0.8410443678739625
$ws.Range("N9").Value = 1.098075133910257
$ws.Range("O9").Value = 1.308740759541166
$ws.Range("B10").Value = 0.1815862360846978
$ws.Range("D10").Value = 0.02460306634682752
$ws.Range("E10").Value = 0.6418191604632568
$ws.Range("F10").Value = 0.440750193474031
$ws.Range("G10").Value = 0.2852451615439975
$ws.Range("H10").Value = 0.4180122960958954
$ws.Range("K10").Value = 0.9879655934633718
$ws.Range("N10").Value = 1.086387834586539
$ws.Range("O10").Value = 1.329805625245228
$ws.Range("B11").Value = 0.1899080348073596
$ws.Range("D11").Value = 0.02624953252309581
$ws.Range("E11").Value = 0.6875660724692665
$ws.Range("F11").Value = 0.4471448160782217
$ws.Range("G11").Value = 0.2902552719902189
$ws.Range("H11").Value = 0.4179476561639177
$ws.Range("K11").Value = 1.05453463850435
$ws.Range("N11").Value = 1.081534409573983
$ws.Range("O11").Value = 1.340464715721652
$ws.Range("B12").Value = 0.1930688545902228
$ws.Range("D12").Value = 0.02687190909173864
$ws.Range("E12").Value = 0.7049121733138861
$ws.Range("F12").Value = 0.4496156591469713
$ws.Range("G12").Value = 0.2921912810256515
$ws.Range("H12").Value = 0.4179588178129165
$ws.Range("K12").Value = 1.079703441420634
$ws.Range("N12").Value = 1.079763043806231
$ws.Range("O12").Value = 1.344656568522879
$ws.Range("B13").Value = 0.1923876946591605
$ws.Range("D13").Value = 0.02673791891829325
$ws.Range("E13").Value = 0.7011753422889484
$ws.Range("F13").Value = 0.4490813219313807
$ws.Range("G13").Value = 0.2917725974444068
$ws.Range("H13").Value = 0.4179548280781376
$ws.Range("K13").Value = 1.074284662611831
$ws.Range("N13").Value = 1.080141581888235
$ws.Range("O13").Value = 1.343746851933417
$ws.Range("B14").Value = 0.1901678874207562
$ws.Range("D14").Value = 0.02630075818205313
$ws.Range("E14").Value = 0.6889926813662157
$ws.Range("F14").Value = 0.4473471040703032
$ws.Range("G14").Value = 0.2904137697028517
$ws.Range("H14").Value = 0.4179478598844781
$ws.Range("K14").Value = 1.056606088424076
$ws.Range("N14").Value = 1.081387345538701
$ws.Range("O14").Value = 1.340806461915605
$ws.Range("B15").Value = 0.188809427112318
$ws.Range("D15").Value = 0.02603283959023628
$ws.Range("E15").Value = 0.6815334604373504
$ws.Range("F15").Value = 0.4462912757355397
$ws.Range("G15").Value = 0.2895865078562849
$ws.Range("H15").Value = 0.4179482344588621
$ws.Range("K15").Value = 1.04577227890897
$ws.Range("N15").Value = 1.082159072288412
$ws.Range("O15").Value = 1.339025657152689
$ws.Range("B16").Value = 0.1810437323711511
$ws.Range("D16").Value = 0.02449531353814649
$ws.Range("E16").Value = 0.6388325833788002
$ws.Range("F16").Value = 0.4403391850679057
$ws.Range("G16").Value = 0.2849231524730413
$ws.Range("H16").Value = 0.418021504947518
$ws.Range("K16").Value = 0.9836096933748877
$ws.Range("N16").Value = 1.086714330491944
$ws.Range("O16").Value = 1.329130749868057
$ws.Range("B17").Value = 0.176296931377891
$ws.Range("D17").Value = 0.02355016675299026
$ws.Range("E17").Value = 0.6126756985582489
$ws.Range("F17").Value = 0.4367754720381214
$ws.Range("G17").Value = 0.2821311378158526
$ws.Range("H17").Value = 0.4181298786557477
$ws.Range("K17").Value = 0.9454059528473522
$ws.Range("N17").Value = 1.089627405503116
$ws.Range("O17").Value = 1.323336762712756
$ws.Range("B18").Value = 0.1735730837329044
$ws.Range("D18").Value = 0.02300584899285241
$ws.Range("E18").Value = 0.5976446473135724
$ws.Range("F18").Value = 0.4347578958043812
$ws.Range("G18").Value = 0.2805504439716486
$ws.Range("H18").Value = 0.4182155046495524
$ws.Range("K18").Value = 0.9234071898376328
$ws.Range("N18").Value = 1.091346531279441
$ws.Range("O18").Value = 1.320105506600243
$ws.Range("B19").Value = 0.1726519388449361
$ws.Range("D19").Value = 0.02282143432149297
$ws.Range("E19").Value = 0.5925576937369925
$ws.Range("F19").Value = 0.4340803006156975
$ws.Range("G19").Value = 0.2800195680491839
$ws.Range("H19").Value = 0.4182484951031284
$ws.Range("K19").Value = 0.9159545319227504
$ws.Range("N19").Value = 1.091936089070806
$ws.Range("O19").Value = 1.319028835237702
$ws.Range("B20").Value = 0.1768015769168727
$ws.Range("D20").Value = 0.02365085139229706
$ws.Range("E20").Value = 0.6154587140373877
$ws.Range("F20").Value = 0.4371515039557892
$ws.Range("G20").Value = 0.2824257425932757
$ws.Range("H20").Value = 0.4181159310613651
$ws.Range("K20").Value = 0.9494754003173398
$ws.Range("N20").Value = 1.089312791666941
$ws.Range("O20").Value = 1.323943054289742
$ws.Range("B21").Value = 0.190819641456244
$ws.Range("D21").Value = 0.02642919318054027
$ws.Range("E21").Value = 0.6925703937951226
$ws.Range("F21").Value = 0.4478551455571065
$ws.Range("G21").Value = 0.2908118355028932
$ws.Range("H21").Value = 0.4179489389676974
$ws.Range("K21").Value = 1.061799794208184
$ws.Range("N21").Value = 1.081019629853145
$ws.Range("O21").Value = 1.341665900484912
$ws.Range("B22").Value = 0.2000367408057429
$ws.Range("D22").Value = 0.02823854410111437
$ws.Range("E22").Value = 0.7431006487879017
$ws.Range("F22").Value = 0.4551382515342652
$ws.Range("G22").Value = 0.296518847432381
$ws.Range("H22").Value = 0.4180475572250515
$ws.Range("K22").Value = 1.134979792685044
$ws.Range("N22").Value = 1.075987267571151
$ws.Range("O22").Value = 1.354155414329227
$ws.Range("B23").Value = 0.1951123919757549
$ws.Range("D23").Value = 0.02727346339954551
$ws.Range("E23").Value = 0.7161189506309142
$ws.Range("F23").Value = 0.4512247458090712
$ws.Range("G23").Value = 0.2934521209142957
$ws.Range("H23").Value = 0.4179758963362445
$ws.Range("K23").Value = 1.095943727451072
$ws.Range("N23").Value = 1.078637686604594
$ws.Range("O23").Value = 1.347406347243918
$ws.Range("B24").Value = 0.1765734106068493
$ws.Range("D24").Value = 0.02360533479282623
$ws.Range("E24").Value = 0.6142004914321717
$ws.Range("F24").Value = 0.4369814025952437
$ws.Range("G24").Value = 0.2822924755617464
$ws.Range("H24").Value = 0.4181221641273538
$ws.Range("K24").Value = 0.9476357118150815
$ws.Range("N24").Value = 1.089454890393867
$ws.Range("O24").Value = 1.323668639087742
$ws.Range("B25").Value = 0.15682103730704
$ws.Range("D25").Value = 0.01962226727788874
$ws.Range("E25").Value = 0.5047815588427227
$ws.Range("F25").Value = 0.4228471017036171
$ws.Range("G25").Value = 0.2712172608834891
$ws.Range("H25").Value = 0.4191702351819373
$ws.Range("K25").Value = 0.7867106198928013
$ws.Range("N25").Value = 1.102790379958961
$ws.Range("O25").Value = 1.301933377085987
